$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1714285714285714
$ws.Range("C2").Value = 0.5904761904761905
$ws.Range("J2").Value = 0.009523809523809525
$ws.Range("P2").Value = 0.1142857142857143
$ws.Range("S2").Value = 0.1142857142857143
$ws.Range("C3").Value = 0.04477611940298507
$ws.Range("J3").Value = 0.04477611940298507
$ws.Range("P3").Value = 0.7164179104477612
$ws.Range("S3").Value = 0.1940298507462687
$ws.Range("J4").Value = 0.1578947368421053
$ws.Range("P4").Value = 0.4736842105263158
$ws.Range("S4").Value = 0.3684210526315789
$ws.Range("B6").Value = 0.07207207207207207
$ws.Range("F6").Value = 0.04504504504504504
$ws.Range("J6").Value = 0.2522522522522522
$ws.Range("O6").Value = 0.009009009009009009
$ws.Range("Q6").Value = 0.1531531531531531
$ws.Range("R6").Value = 0.05405405405405406
$ws.Range("S6").Value = 0.4144144144144144
$ws.Range("B7").Value = 0.09375
$ws.Range("F7").Value = 0.0625
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.0625
$ws.Range("Q7").Value = 0.1354166666666667
$ws.Range("R7").Value = 0.05208333333333334
$ws.Range("S7").Value = 0.46875
$ws.Range("B8").Value = 0.0761904761904762
$ws.Range("D8").Value = 0.01904761904761905
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.1476190476190476
$ws.Range("O8").Value = 0.01904761904761905
$ws.Range("Q8").Value = 0.1619047619047619
$ws.Range("R8").Value = 0.1380952380952381
$ws.Range("S8").Value = 0.3714285714285714
$ws.Range("B9").Value = 0.07954545454545454
$ws.Range("D9").Value = 0.01136363636363636
$ws.Range("F9").Value = 0.03409090909090909
$ws.Range("J9").Value = 0.1704545454545454
$ws.Range("O9").Value = 0.01136363636363636
$ws.Range("Q9").Value = 0.1931818181818182
$ws.Range("R9").Value = 0.04545454545454546
$ws.Range("S9").Value = 0.4545454545454545
$ws.Range("B10").Value = 0.06241519674355495
$ws.Range("D10").Value = 0.0203527815468114
$ws.Range("E10").Value = 0.00135685210312076
$ws.Range("F10").Value = 0.06648575305291723
$ws.Range("J10").Value = 0.1248303934871099
$ws.Range("O10").Value = 0.01085481682496608
$ws.Range("Q10").Value = 0.2116689280868385
$ws.Range("R10").Value = 0.07734056987788331
$ws.Range("S10").Value = 0.4246947082767978
$ws.Range("G11").Value = 0.1257142857142857
$ws.Range("J11").Value = 0.16
$ws.Range("K11").Value = 0.24
$ws.Range("L11").Value = 0.4628571428571429
$ws.Range("S11").Value = 0.01142857142857143
$ws.Range("G12").Value = 0.6463414634146342
$ws.Range("J12").Value = 0.3048780487804878
$ws.Range("K12").Value = 0.01219512195121951
$ws.Range("L12").Value = 0.01219512195121951
$ws.Range("S12").Value = 0.02439024390243903
$ws.Range("G13").Value = 0.7241379310344828
$ws.Range("J13").Value = 0.2758620689655172
$ws.Range("F15").Value = 0.02061855670103093
$ws.Range("H15").Value = 0.154639175257732
$ws.Range("I15").Value = 0.09278350515463918
$ws.Range("J15").Value = 0.3711340206185567
$ws.Range("K15").Value = 0.06185567010309279
$ws.Range("M15").Value = 0.01030927835051546
$ws.Range("O15").Value = 0.03092783505154639
$ws.Range("S15").Value = 0.2577319587628866
$ws.Range("F16").Value = 0.01449275362318841
$ws.Range("H16").Value = 0.1449275362318841
$ws.Range("I16").Value = 0.05797101449275362
$ws.Range("J16").Value = 0.4927536231884058
$ws.Range("K16").Value = 0.1014492753623188
$ws.Range("M16").Value = 0.01449275362318841
$ws.Range("O16").Value = 0.08695652173913043
$ws.Range("S16").Value = 0.08695652173913043
$ws.Range("H17").Value = 0.1324786324786325
$ws.Range("I17").Value = 0.1239316239316239
$ws.Range("J17").Value = 0.5213675213675214
$ws.Range("K17").Value = 0.07264957264957266
$ws.Range("M17").Value = 0.0170940170940171
$ws.Range("O17").Value = 0.04700854700854701
$ws.Range("S17").Value = 0.08547008547008547
$ws.Range("F18").Value = 0.02
$ws.Range("H18").Value = 0.16
$ws.Range("I18").Value = 0.08
$ws.Range("J18").Value = 0.53
$ws.Range("K18").Value = 0.09
$ws.Range("O18").Value = 0.06
$ws.Range("S18").Value = 0.06
$ws.Range("F19").Value = 0.02074074074074074
$ws.Range("H19").Value = 0.2074074074074074
$ws.Range("I19").Value = 0.0562962962962963
$ws.Range("J19").Value = 0.3851851851851852
$ws.Range("K19").Value = 0.1377777777777778
$ws.Range("M19").Value = 0.03555555555555556
$ws.Range("O19").Value = 0.0562962962962963
$ws.Range("S19").Value = 0.1007407407407407
